$wb = $excel.ActiveWorkbook

# Insert a new worksheet before the first sheet (AddCustomerTest) to hold the test suite run-mode table.
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "test_suite"

# Populate the run-mode table.
$newSheet.Range("A1").Value = "TCID"
$newSheet.Range("B1").Value = "Runmode"
$newSheet.Range("A2").Value = "BankManagerLoginTest"
$newSheet.Range("B2").Value = "Y"
$newSheet.Range("A3").Value = "AddCustomerTest"
$newSheet.Range("B3").Value = "Y"
$newSheet.Range("A4").Value = "OpenAccountTest"
$newSheet.Range("B4").Value = "N"

# Match column sizing / selection seen in the authored workbook.
$newSheet.Columns.Item(1).ColumnWidth = 21
$newSheet.Range("B4").Select() | Out-Null
